# edit.ps1
# Applies the tts.pptx diff:
#  1. Slide 2 (index 2): add a new title placeholder shape with text "TTS"
#     positioned/sized per the diff.
#  2. Slide 2 (index 2): add a new text box shape with text
#     "remove speech background sound(music/noise and so on)" positioned/
#     sized per the diff.
#  3. Slide 4 (index 4): change the title text from "nlp" to "NLP".
#
# Shape.Left/Top/Width/Height (and AddTextbox's coordinate args) are stored
# internally as single-precision (float32) point values and then converted
# to EMU by flooring (pt * 12700). A plain "emu / 12700.0" division can,
# after the float64->float32 demotion, land a hair below the intended value
# and truncate to one EMU less than the target (e.g. 1077595 EMU becomes
# 1077594). EmuToPt nudges the point value up by tiny increments until the
# float32 round-trip reproduces the exact requested EMU, so the written
# OOXML matches the target byte-for-byte on geometry.

function EmuToPt {
    param([double]$emu)
    $pt = $emu / 12700.0
    for ($i = 0; $i -lt 4000; $i++) {
        $trial = $pt + ($i * 0.0000001)
        $f32 = [float]$trial
        $emuTrial = [Math]::Floor([double]$f32 * 12700.0)
        if ($emuTrial -eq $emu) {
            return $trial
        }
    }
    return $pt
}

$p = $ppt.ActivePresentation
$slide2 = $p.Slides.Item(2)
$slide4 = $p.Slides.Item(4)

# Locate slide 4's title placeholder robustly (falls back to shape 1, which
# is the title on this slide).
$slide4Title = $slide4.Shapes.Item(1)
for ($i = 1; $i -le $slide4.Shapes.Placeholders.Count; $i++) {
    $ph = $slide4.Shapes.Placeholders.Item($i)
    if ($ph.PlaceholderFormat.Type -eq 1) {
        $slide4Title = $ph
    }
}

# ---------------------------------------------------------------------
# 1. New title shape ("TTS") on slide 2.
#    Slide 2's layout already defines a title placeholder, but the slide
#    itself has none yet. Copy the existing title placeholder shape from
#    slide 4 (which already carries the correct <p:ph type="title"/> +
#    style wiring) and paste it onto slide 2, then retarget its name,
#    geometry and text.
# ---------------------------------------------------------------------
$slide4Title.Copy()
$titlePasted = $slide2.Shapes.Paste()
$titleShape = $titlePasted.Item(1)

$titleShape.Name = "标题 45"
$titleShape.TextFrame.TextRange.Text = "TTS"
$titleShape.Left = EmuToPt 487680
$titleShape.Top = EmuToPt 95885
$titleShape.Width = EmuToPt 1449705
$titleShape.Height = EmuToPt 1077595

# ---------------------------------------------------------------------
# 2. New text box ("remove speech background sound(music/noise and so
#    on)") on slide 2. Copy an existing plain text box on the same slide
#    (same bodyPr/run formatting as the target) so the generated XML
#    attributes line up exactly, then retarget name, geometry and text.
# ---------------------------------------------------------------------
$templateShape = $null
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "TTS main framework") {
                $templateShape = $shp
            }
        }
    }
}

$templateShape.Copy()
$boxPasted = $slide2.Shapes.Paste()
$textBoxShape = $boxPasted.Item(1)

$textBoxShape.Name = "文本框 46"
$textBoxShape.TextFrame.TextRange.Text = "remove speech background sound(music/noise and so on)"
$textBoxShape.Left = EmuToPt 4519295
$textBoxShape.Top = EmuToPt 6489700
$textBoxShape.Width = EmuToPt 5581650
$textBoxShape.Height = EmuToPt 368300

# ---------------------------------------------------------------------
# 3. Fix the "nlp" -> "NLP" typo on slide 4's title.
# ---------------------------------------------------------------------
$slide4Title.TextFrame.TextRange.Text = "NLP"
